# Apply "Added headers and source link indication" edit.
#
# 1) Bold the header row (A1:B1 -> "Title"/"URL") on both worksheets.
# 2) Bold the "source link" indicator cells in column C of "Second run".
# 3) Remove the now-obsolete crawl rows (Review / essdd-5-221 / essdd-5-243 /
#    essd-4-75 / pr_copernicus_article_level_metrics / review_process...)
#    from the "Second run" sheet - deleting these rows also drops their
#    now-unused shared strings automatically.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "First run"
$ws2 = $wb.Worksheets.Item(2)   # "Second run"

# --- 1) Bold header cells on both sheets --------------------------------
$ws1.Range("A1:B1").Font.Bold = $true
$ws2.Range("A1:B1").Font.Bold = $true

# --- 2) Bold every "source link" cell (column C) on the second sheet ----
$sourceLinkRows = @(2,13,60,77,110,112,120,125,139,151,159,179,182,187,195)
foreach ($r in $sourceLinkRows) {
    $ws2.Range("C" + $r).Font.Bold = $true
}

# --- 3) Delete the rows whose content is being dropped -------------------
# Delete from the bottom up so earlier row numbers stay valid.
$rowsToDelete = @(53,35,34,33,22)
foreach ($r in $rowsToDelete) {
    $ws2.Rows.Item($r).Delete()
}
